$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$newDate = [DateTime]::FromOADate(45184)

for ($r = 2; $r -le 115; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
